$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row: "_old" -> "_FV2304" and "_new" -> "_FV2310" ---
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_old$", "_FV2304")
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_new$", "_FV2310")
}

# --- Freeze top row (pane split after row 1) ---
$ws.Activate()
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Create table (ListObject) over the used range A1:U63 ---
# 1 = xlSrcRange (source type), 1 = xlYes (range has headers)
$tableRange = $ws.Range("A1:U63")
$list = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$list.Name = "Table1"
